# Insert a new record row right before the existing row 545. This shifts the
# previous rows 545:634 down to 546:635 (preserving all their data untouched)
# and the sheet's used range grows from A1:T634 to A1:T635.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(545).Insert()

# Populate the newly inserted row 545 with the new "Uva" price record.
$ws.Range("A545").Value = 9
$ws.Range("B545").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C545").Value = "Metropolitana"
$ws.Range("D545").Value = 44711
$ws.Range("E545").Value = 13
$ws.Range("F545").Value = "Fruta"
$ws.Range("G545").Value = 100109
$ws.Range("H545").Value = "Uva"
$ws.Range("I545").Value = 100109001
$ws.Range("J545").Value = "Uva"
$ws.Range("K545").Value = "Red Globe"
$ws.Range("L545").Value = "Primera"
$ws.Range("M545").Value = 95
$ws.Range("N545").Value = 9000
$ws.Range("O545").Value = 9000
$ws.Range("P545").Value = 9000
$ws.Range("Q545").Value = "$/bandeja 18 kilos"
$ws.Range("R545").Value = "Región de O'Higgins"
$ws.Range("S545").Value = 500
$ws.Range("T545").Value = 18
